$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'63.356.57"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -0.35%  '

$ws.Range('D3').Value = "'3.190.19"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -3.56%  '

$ws.Range('D4').Value = "'0.999"
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.13%  '

$ws.Range('D5').Value = "'593.80"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.68%  '

$ws.Range('D6').Value = "'136.24"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -3.23%  '

$ws.Range('E7').Value = '  -0.10%  '

$ws.Range('D8').Value = "'3.189.20"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -3.58%  '

$ws.Range('E9').Value = '  -0.50%  '

$ws.Range('D10').Value = "'0.142"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -4.66%  '

$ws.Range('D11').Value = "'5.27"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -4.48%  '

$ws.Range('D12').Value = "'0.456"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -2.32%  '

$ws.Range('D13').Value = "'0.0000239"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -2.51%  '

$ws.Range('D14').Value = "'34.70"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.15%  '

$ws.Range('D15').Value = "'3.714.85"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -3.58%  '

$ws.Range('E16').Value = '  -1.72%  '

$ws.Range('D17').Value = "'3.183.94"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -3.69%  '

$ws.Range('D18').Value = "'63.282.07"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.68%  '

$ws.Range('D19').Value = "'6.59"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -3.69%  '

$ws.Range('D20').Value = "'463.12"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -3.23%  '

$ws.Range('D21').Value = "'14.01"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.23%  '

$ws.Range('D22').Value = "'0.701"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -4.76%  '

$ws.Range('D23').Value = "'7.66"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -3.83%  '

$ws.Range('D24').Value = "'13.42"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -3.41%  '

$ws.Range('D25').Value = "'82.71"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -2.94%  '

$ws.Range('E26').Value = '  +0.15%  '

$ws.Range('B27').Value = 'PancakeSwap'
$ws.Range('C27').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D27').Value = "'2.68"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -3.25%  '

$ws.Range('B28').Value = 'FirstDigitalUSD'
$ws.Range('C28').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D28').Value = "'0.998"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.27%  '

$ws.Range('D29').Value = "'7.76"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -4.56%  '

$ws.Range('D30').Value = "'6.80"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -5.15%  '

$ws.Range('D31').Value = "'2.05"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -4.05%  '

$ws.Range('D32').Value = "'27.42"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -4.24%  '

$ws.Range('D33').Value = "'0.103"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -2.27%  '

$ws.Range('D34').Value = "'2.39"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -4.86%  '

$ws.Range('E35').Value = '  -5.32%  '

$ws.Range('D36').Value = "'5.85"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -2.98%  '

$ws.Range('D37').Value = "'51.35"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -2.06%  '

$ws.Range('D38').Value = "'0.0₃0721"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -1.73%  '

$ws.Range('D39').Value = "'0.0388"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -2.69%  '

$ws.Range('B40').Value = 'dogwifhat'
$ws.Range('C40').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D40').Value = "'2.70"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.41%  '

$ws.Range('B41').Value = 'Bittensor'
$ws.Range('C41').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D41').Value = "'403.20"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -6.01%  '

$ws.Range('B42').Value = 'Cosmos'
$ws.Range('C42').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D42').Value = "'8.12"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -2.15%  '

$ws.Range('D43').Value = "'0.112"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -6.22%  '

$ws.Range('D44').Value = "'2.812.35"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -10.07%  '

$ws.Range('D45').Value = "'0.254"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -3.51%  '

$ws.Range('B46').Value = 'Fetch.AI'
$ws.Range('C46').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D46').Value = "'2.14"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -3.18%  '

$ws.Range('B47').Value = 'USDe'
$ws.Range('C47').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D47').Value = "'0.999"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.03%  '

$ws.Range('D48').Value = "'127.03"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.88%  '

$ws.Range('B49').Value = 'Arweave'
$ws.Range('C49').Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range('D49').Value = "'35.51"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -2.71%  '

$ws.Range('B50').Value = 'InjectiveProtocol'
$ws.Range('C50').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D50').Value = "'25.39"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -3.31%  '

$ws.Range('D51').Value = "'0.112"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -1.56%  '
